$wb = $excel.ActiveWorkbook

# --- Update the Date metadata value on the "Metadata" sheet ---
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2026-01-23T08:28:04+00:00"

# --- Fix the source mapping on "Mapping Table 1" row 5: ---
# FRCDADICOMExamenImagerie.description -> FRCDADICOMExamenImagerie.text
$mapSheet1 = $wb.Worksheets.Item("Mapping Table 1")
$mapSheet1.Range("A5").Value = "FRCDADICOMExamenImagerie.text"

$wb.Save()
